# "Fruta / hortaliza, semanal"
#
# The weekly refresh inserts one new observation (row) at the top of the
# data table (row 37, right after the first 35 weekly rows) and pushes
# every existing record below it down by one row. This script reproduces
# that by inserting a new row at row 37 and filling it with the new
# weekly price record; Excel automatically shifts all the rows that used
# to be at 37..94 down to 38..95 (carrying their values/styles with them),
# which is exactly what the target workbook looks like.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37, pushing old rows 37-94 to 38-95.
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the new weekly record.
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44894
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112040
$ws.Range("G37").Value = "Cilantro"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 250
$ws.Range("K37").Value = 1500
$ws.Range("L37").Value = 1800
$ws.Range("M37").Value = 1620
$ws.Range("N37").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 810
$ws.Range("Q37").Value = 2
$ws.Range("R37").Value = "Hortaliza"

# Keep the same date formatting the rest of the date column uses.
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
